$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for rows 2-8 (repulled data / mean calculation)
$ws.Range("F2").Value = -3
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = 7
$ws.Range("F8").Value = 3
